$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinate values to whole numbers
$ws.Range("Q2").Value = 541566
$ws.Range("R2").Value = 7246381
$ws.Range("Q3").Value = 541566
$ws.Range("R3").Value = 7246381

# Clear the start/end time cells (Z and AB columns) for rows 2 and 3
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
